$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 7-26 per diff (only changed cells) ---
# Row 7
$ws.Range("D7").Value = 44525
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L7").Value = 'Especial'
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 25200
$ws.Range("O7").Value = 25200
$ws.Range("P7").Value = 25200
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("S7").Value = 1400
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44525
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 21600
$ws.Range("O8").Value = 21600
$ws.Range("P8").Value = 21600
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("S8").Value = 1200
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44159
$ws.Range("D9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K9").Value = 'Castle Brite'
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("Q9").Value = '$/bandeja 10 kilos'
$ws.Range("S9").Value = 800
$ws.Range("T9").Value = 10

# Row 10
$ws.Range("D10").Value = 44159
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K10").Value = 'Castle Brite'
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 65
$ws.Range("N10").Value = 7000
$ws.Range("O10").Value = 7000
$ws.Range("P10").Value = 7000
$ws.Range("R10").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S10").Value = 700

# Row 11
$ws.Range("D11").Value = 44519
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 22500
$ws.Range("O11").Value = 22500
$ws.Range("P11").Value = 22500
$ws.Range("Q11").Value = '$/caja 15 kilos granel'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1500
$ws.Range("T11").Value = 15

# Row 12
$ws.Range("D12").Value = 44175
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K12").Value = 'Modesto'
$ws.Range("M12").Value = 140
$ws.Range("N12").Value = 11000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 11571
$ws.Range("Q12").Value = '$/caja 12 kilos'
$ws.Range("S12").Value = 964
$ws.Range("T12").Value = 12

# Row 13
$ws.Range("D13").Value = 44168
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K13").Value = 'Dina'
$ws.Range("M13").Value = 40
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 14000
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 1400

# Row 14
$ws.Range("D14").Value = 44162
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M14").Value = 70
$ws.Range("N14").Value = 8500
$ws.Range("O14").Value = 8500
$ws.Range("P14").Value = 8500
$ws.Range("S14").Value = 850

# Row 15
$ws.Range("D15").Value = 44162
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M15").Value = 75
$ws.Range("N15").Value = 14000
$ws.Range("P15").Value = 14400
$ws.Range("S15").Value = 800

# Row 16
$ws.Range("D16").Value = 44167
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K16").Value = 'Castle Brite'
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 85
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 10000
$ws.Range("R16").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S16").Value = 1000

# Row 17
$ws.Range("D17").Value = 44167
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K17").Value = 'Castle Brite'
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 9500
$ws.Range("O17").Value = 9500
$ws.Range("P17").Value = 9500
$ws.Range("Q17").Value = '$/bandeja 10 kilos'
$ws.Range("R17").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S17").Value = 950
$ws.Range("T17").Value = 10

# Row 18
$ws.Range("D18").Value = 44167
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K18").Value = 'Castle Brite'
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 15000
$ws.Range("R18").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S18").Value = 833

# Row 19
$ws.Range("D19").Value = 44174
$ws.Range("D19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K19").Value = 'Modesto'
$ws.Range("M19").Value = 120
$ws.Range("N19").Value = 8500
$ws.Range("O19").Value = 8500
$ws.Range("P19").Value = 8500
$ws.Range("Q19").Value = '$/bandeja 10 kilos'
$ws.Range("R19").Value = 'Región Metropolitana'
$ws.Range("S19").Value = 850
$ws.Range("T19").Value = 10

# Row 20
$ws.Range("D20").Value = 44174
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K20").Value = 'Modesto'
$ws.Range("M20").Value = 180
$ws.Range("N20").Value = 15000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 15000
$ws.Range("Q20").Value = '$/caja 18 kilos'
$ws.Range("R20").Value = 'Región Metropolitana'
$ws.Range("S20").Value = 833
$ws.Range("T20").Value = 18

# Row 21
$ws.Range("D21").Value = 44174
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K21").Value = 'Modesto'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 120
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 12000
$ws.Range("Q21").Value = '$/caja 18 kilos'
$ws.Range("R21").Value = 'Región Metropolitana'
$ws.Range("S21").Value = 667
$ws.Range("T21").Value = 18

# Row 22
$ws.Range("D22").Value = 44189
$ws.Range("D22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M22").Value = 130
$ws.Range("N22").Value = 12000
$ws.Range("O22").Value = 12000
$ws.Range("P22").Value = 12000
$ws.Range("Q22").Value = '$/caja 18 kilos'
$ws.Range("S22").Value = 667
$ws.Range("T22").Value = 18

# Row 23
$ws.Range("D23").Value = 44523
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K23").Value = 'Castle Brite'
$ws.Range("M23").Value = 320
$ws.Range("N23").Value = 10000
$ws.Range("O23").Value = 10000
$ws.Range("P23").Value = 10000
$ws.Range("Q23").Value = '$/bandeja 10 kilos'
$ws.Range("S23").Value = 1000
$ws.Range("T23").Value = 10

# Row 24
$ws.Range("D24").Value = 44169
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K24").Value = 'Dina'
$ws.Range("M24").Value = 80
$ws.Range("R24").Value = 'Región de O''Higgins'

# Row 25
$ws.Range("D25").Value = 44195
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K25").Value = 'Patterson'
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 124
$ws.Range("N25").Value = 13000
$ws.Range("O25").Value = 13000
$ws.Range("P25").Value = 13000
$ws.Range("Q25").Value = '$/caja 15 kilos'
$ws.Range("R25").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S25").Value = 867
$ws.Range("T25").Value = 15

# Row 26
$ws.Range("D26").Value = 44176
$ws.Range("D26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K26").Value = 'Modesto'
$ws.Range("M26").Value = 115
$ws.Range("N26").Value = 11000
$ws.Range("O26").Value = 12000
$ws.Range("P26").Value = 11609
$ws.Range("Q26").Value = '$/caja 12 kilos'
$ws.Range("S26").Value = 967
$ws.Range("T26").Value = 12

# --- Add new rows 27-29 ---
# Row 27
$ws.Range("A27").Value = 9
$ws.Range("B27").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C27").Value = 'Metropolitana'
$ws.Range("D27").Value = 44166
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 'Fruta'
$ws.Range("G27").Value = 100103
$ws.Range("H27").Value = 'Frutos de hueso (carozo)'
$ws.Range("I27").Value = 100103003
$ws.Range("J27").Value = 'Damasco'
$ws.Range("K27").Value = 'Castle Brite'
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 120
$ws.Range("N27").Value = 10000
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 10000
$ws.Range("Q27").Value = '$/bandeja 10 kilos'
$ws.Range("R27").Value = 'Región Metropolitana'
$ws.Range("S27").Value = 1000
$ws.Range("T27").Value = 10

# Row 28
$ws.Range("A28").Value = 9
$ws.Range("B28").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C28").Value = 'Metropolitana'
$ws.Range("D28").Value = 44166
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 'Fruta'
$ws.Range("G28").Value = 100103
$ws.Range("H28").Value = 'Frutos de hueso (carozo)'
$ws.Range("I28").Value = 100103003
$ws.Range("J28").Value = 'Damasco'
$ws.Range("K28").Value = 'Castle Brite'
$ws.Range("L28").Value = 'Segunda'
$ws.Range("M28").Value = 120
$ws.Range("N28").Value = 8000
$ws.Range("O28").Value = 8000
$ws.Range("P28").Value = 8000
$ws.Range("Q28").Value = '$/bandeja 10 kilos'
$ws.Range("R28").Value = 'Región Metropolitana'
$ws.Range("S28").Value = 800
$ws.Range("T28").Value = 10

# Row 29
$ws.Range("A29").Value = 9
$ws.Range("B29").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C29").Value = 'Metropolitana'
$ws.Range("D29").Value = 44194
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 'Fruta'
$ws.Range("G29").Value = 100103
$ws.Range("H29").Value = 'Frutos de hueso (carozo)'
$ws.Range("I29").Value = 100103003
$ws.Range("J29").Value = 'Damasco'
$ws.Range("K29").Value = 'Patterson'
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 120
$ws.Range("N29").Value = 13000
$ws.Range("O29").Value = 13000
$ws.Range("P29").Value = 13000
$ws.Range("Q29").Value = '$/caja 15 kilos'
$ws.Range("R29").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S29").Value = 867
$ws.Range("T29").Value = 15

Write-Host "Edit complete."
